# Resumen Downtime - update report with Home+ASH breakdown
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column width: A from 30 -> 35 characters -------------------------------
# Excel's ColumnWidth property and the stored OOXML <col width> differ by the
# standard 5px/MaximumDigitWidth padding offset (~0.8333 for Calibri 11).
$ws.Columns("A").ColumnWidth = 35 - 5/6

# --- Plain text label/value updates (not numeric-looking, safe as .Value) --
$ws.Range("A2").Value = "Uptime/Downtime (% General)"
$ws.Range("A3").Value = "Uptime/Downtime (Home+ASH hh:mm:ss)"
$ws.Range("A6").Value = "URLs analizadas:"
$ws.Range("B6").Value = "Home (main) + ASH"
$ws.Range("A7").Value = "Downtime Home (segundos):"
$ws.Range("A8").Value = "Downtime ASH (segundos):"
$ws.Range("A9").Value = "Total downtime (segundos):"
$ws.Range("A10").Value = "Tiempo total sistema (segundos):"
$ws.Range("A11").Value = "Fecha de generación:"

# --- Numeric-looking values must stay stored as literal text (inline/shared
# string), matching the source report. Use an off-sheet scratch cell forced
# to Text format, then paste only the *value* into the destination so the
# destination's own cell style/format is left completely untouched. -------
$scratch = $ws.Range("D1")
$scratch.NumberFormat = "@"

$scratch.Value = "0.0000%"
$scratch.Copy()
$ws.Range("B2").PasteSpecial(-4163)

$scratch.Value = "00:00:00"
$scratch.Copy()
$ws.Range("B3").PasteSpecial(-4163)

$scratch.Value = "0"
$scratch.Copy()
$ws.Range("B7").PasteSpecial(-4163)

$scratch.Value = "0"
$scratch.Copy()
$ws.Range("B8").PasteSpecial(-4163)

$scratch.Value = "0"
$scratch.Copy()
$ws.Range("B9").PasteSpecial(-4163)

$scratch.Value = "2418547"
$scratch.Copy()
$ws.Range("B10").PasteSpecial(-4163)

$scratch.Value = "2025-06-27 14:19:46"
$scratch.Copy()
$ws.Range("B11").PasteSpecial(-4163)

$scratch.Clear()

# --- Style fix-ups -----------------------------------------------------------
# B2 must share the same style as B3 (bold, bordered, centered value cell).
$ws.Range("B3").Copy()
$ws.Range("B2").PasteSpecial(-4122)

# Newly created rows (8, 10, 11) need the plain bordered "info" style that is
# already used throughout rows 6-9 of that block; A6 still carries it.
$ws.Range("A6").Copy()
$ws.Range("A7:B11").PasteSpecial(-4122)

$excel.CutCopyMode = 0
